$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9665321707725525
$ws.Range("B1").Value = 2.236487150192261
$ws.Range("C1").Value = 8.13615894317627
$ws.Range("D1").Value = 1.806392073631287
$ws.Range("E1").Value = 1.224856972694397
